# Generate Report for Handback
#
# The handback process completed for both locales (zh-cn, de-de). This:
#   1) flips the status label from "Ready for handoff" to
#      "Handed back: in sync with en-US" everywhere it is shown
#      (Overview summary sheet + the two per-locale detail sheets),
#   2) fills in the "Latest Target File" / "Latest Handback File" hyperlink
#      columns (E/F) for the two real rows on each locale sheet, and
#   3) stamps the "Latest Handback DateTime" column (G) with the real
#      handback timestamp instead of the epoch placeholder.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1) Overview sheet: just the status label changes (B2:C3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2) Per-locale detail sheets
# ---------------------------------------------------------------------

# zh-cn
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B2").Value = $newStatus
$ws.Range("B3").Value = $newStatus

$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/54405aeff2c9f55c5b8f58c964a65dffbdd73257/e2e/07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md", "", "", "07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c7de5afd40b013df8b5200fa73b16cf5ee4fe089/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/hb/07f22b3b-2bf3-4ca5-b902-c4c6baae3194.14254085f8f7cae4224c47ac3ba3ec67fd5a93b4.zh-cn.xlf", "", "", "07f22b3b-2bf3-4ca5-b902-c4c6baae3194.14254085f8f7cae4224c47ac3ba3ec67fd5a93b4.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/54405aeff2c9f55c5b8f58c964a65dffbdd73257/e2e/ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md", "", "", "ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c7de5afd40b013df8b5200fa73b16cf5ee4fe089/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/hb/ce44ac28-9d01-4602-b293-d67cbbaf5ed3.7214772345b417cebf6937dce591607708667748.zh-cn.xlf", "", "", "ce44ac28-9d01-4602-b293-d67cbbaf5ed3.7214772345b417cebf6937dce591607708667748.zh-cn.xlf")

$ws.Range("G2").Value = "2016-02-17 04:54:49"
$ws.Range("G3").Value = "2016-02-17 04:54:49"

# de-de
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B2").Value = $newStatus
$ws.Range("B3").Value = $newStatus

$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/54405aeff2c9f55c5b8f58c964a65dffbdd73257/e2e/07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md", "", "", "07f22b3b-2bf3-4ca5-b902-c4c6baae3194.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2e73d44a9854f957f9ed515654c8d1c33fa676b7/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/hb/07f22b3b-2bf3-4ca5-b902-c4c6baae3194.14254085f8f7cae4224c47ac3ba3ec67fd5a93b4.de-de.xlf", "", "", "07f22b3b-2bf3-4ca5-b902-c4c6baae3194.14254085f8f7cae4224c47ac3ba3ec67fd5a93b4.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/54405aeff2c9f55c5b8f58c964a65dffbdd73257/e2e/ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md", "", "", "ce44ac28-9d01-4602-b293-d67cbbaf5ed3.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2e73d44a9854f957f9ed515654c8d1c33fa676b7/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/hb/ce44ac28-9d01-4602-b293-d67cbbaf5ed3.7214772345b417cebf6937dce591607708667748.de-de.xlf", "", "", "ce44ac28-9d01-4602-b293-d67cbbaf5ed3.7214772345b417cebf6937dce591607708667748.de-de.xlf")

$ws.Range("G2").Value = "2016-02-17 04:55:07"
$ws.Range("G3").Value = "2016-02-17 04:55:07"
